$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '33.934.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.773.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.29%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.545'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.09%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.01'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.284'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0697'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0921'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.036.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.789.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.930.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.618'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.48'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0764'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.22%  '
$ws.Range("E24").Value = '  -1.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.17'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.10%  '
$ws.Range("E28").Value = '  -3.08%  '
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0513'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.65%  '
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.48'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.64%  '
$ws.Range("E34").Value = '  -4.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.389.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.621'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0184'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.42%  '
$ws.Range("E39").Value = '  +1.75%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.48'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("E44").Value = '  -3.85%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.931.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '103.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.996'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.61'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0121'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.33%  '
